$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay plain text even though the new value looks
    # like a number (e.g. "308.35"), matching the source data exactly
    # instead of letting Excel auto-convert it to a numeric cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2
$ws.Range("D2").Value = "46.425.09"
$ws.Range("E2").Value = "  +2.22%  "

# Row 3
$ws.Range("D3").Value = "2.611.25"
$ws.Range("E3").Value = "  +4.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
Set-TextValue "D5" "308.35"
$ws.Range("E5").Value = "  +3.80%  "

# Row 6
Set-TextValue "D6" "100.36"
$ws.Range("E6").Value = "  +4.01%  "

# Row 7
$ws.Range("E7").Value = "  +2.76%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
Set-TextValue "D9" "0.580"
$ws.Range("E9").Value = "  +7.42%  "

# Row 10
Set-TextValue "D10" "39.56"
$ws.Range("E10").Value = "  +8.31%  "

# Row 11
Set-TextValue "D11" "0.0846"
$ws.Range("E11").Value = "  +5.84%  "

# Row 12
Set-TextValue "D12" "54.24"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
Set-TextValue "D13" "8.16"
$ws.Range("E13").Value = "  +7.51%  "

# Row 14
$ws.Range("D14").Value = "3.012.23"
$ws.Range("E14").Value = "  +3.52%  "

# Row 15
$ws.Range("E15").Value = "  +1.53%  "

# Row 16
$ws.Range("D16").Value = "2.619.16"
$ws.Range("E16").Value = "  +3.51%  "

# Row 17
Set-TextValue "D17" "0.921"
$ws.Range("E17").Value = "  +5.66%  "

# Row 18
Set-TextValue "D18" "14.94"
$ws.Range("E18").Value = "  +3.37%  "

# Row 19
$ws.Range("D19").Value = "46.544.66"
$ws.Range("E19").Value = "  +2.01%  "

# Row 20
$ws.Range("E20").Value = "  +4.89%  "

# Row 21
$ws.Range("E21").Value = "  -1.38%  "

# Row 22
Set-TextValue "D22" "6.74"
$ws.Range("E22").Value = "  +3.84%  "

# Row 23
Set-TextValue "D23" "71.47"
$ws.Range("E23").Value = "  +4.32%  "

# Row 24
Set-TextValue "D24" "273.59"
$ws.Range("E24").Value = "  +9.73%  "

# Row 25
$ws.Range("E25").Value = "  +6.52%  "

# Row 26
$ws.Range("E26").Value = "  +6.29%  "

# Row 27
Set-TextValue "D27" "29.06"
$ws.Range("E27").Value = "  +26.97%  "

# Row 28
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("E29").Value = "  -0.61%  "

# Row 30
Set-TextValue "D30" "10.61"
$ws.Range("E30").Value = "  +5.91%  "

# Row 31
Set-TextValue "D31" "39.03"
$ws.Range("E31").Value = "  -2.44%  "

# Row 32
Set-TextValue "D32" "2.22"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33
Set-TextValue "D33" "6.35"
$ws.Range("E33").Value = "  +10.82%  "

# Row 34
Set-TextValue "D34" "3.65"
$ws.Range("E34").Value = "  -3.41%  "

# Row 35 (ARBITRUM -> WEMIXToken; ranks swapped with row 36)
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D35" "2.86"
$ws.Range("E35").Value = "  +1.80%  "

# Row 36 (WEMIXToken -> ARBITRUM)
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D36" "2.24"
$ws.Range("E36").Value = "  +4.36%  "

# Row 37
Set-TextValue "D37" "0.0840"
$ws.Range("E37").Value = "  +4.77%  "

# Row 38
Set-TextValue "D38" "151.09"
$ws.Range("E38").Value = "  +1.67%  "

# Row 39
$ws.Range("E39").Value = "  +5.38%  "

# Row 40
$ws.Range("E40").Value = "  +4.60%  "

# Row 41
Set-TextValue "D41" "23.26"
$ws.Range("E41").Value = "  +40.63%  "

# Row 42
$ws.Range("E42").Value = "  +2.65%  "

# Row 43
Set-TextValue "D43" "3.66"
$ws.Range("E43").Value = "  +9.21%  "

# Row 44
$ws.Range("E44").Value = "  +7.65%  "

# Row 45
Set-TextValue "D45" "4.10"
$ws.Range("E45").Value = "  +1.12%  "

# Row 46
$ws.Range("D46").Value = "2.125.56"
$ws.Range("E46").Value = "  +5.92%  "

# Row 47
Set-TextValue "D47" "0.998"
$ws.Range("E47").Value = "  -0.45%  "

# Row 48
Set-TextValue "D48" "93.81"
$ws.Range("E48").Value = "  +3.74%  "

# Row 49
Set-TextValue "D49" "9.51"
$ws.Range("E49").Value = "  +8.54%  "

# Row 50
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
Set-TextValue "D51" "109.22"
$ws.Range("E51").Value = "  +3.82%  "
